$d = $word.ActiveDocument

# The document currently ends with an empty paragraph (right before the
# section properties). Add a new (also initially empty) paragraph after it,
# then fill that new paragraph - plus one more right after - via InsertXML so
# we get exact control over run/paragraph-mark formatting (bold + complex
# script bold + single underline on the heading, highlighted run on the
# command line), matching what "git clone" documentation entry looks like.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$fragment = @"
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>clonar repositorio</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>git</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> clone </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>https del repositorio</w:t>
  </w:r>
</w:p>
"@

$newPara.Range.InsertXML($fragment)
